$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the previous scrape's long tail (old rows 12-32); remaining rows shift up.
#    New used range / dimension becomes A1:H11.
$ws.Rows("12:32").Delete()

# 2) Narrow columns B and H slightly (target raw stored width 50 / 12).
#    ColumnWidth is run through Excel's MDW pixel-rounding before being stored,
#    so the inputs below are chosen to land exactly on the desired stored width.
$ws.Columns.Item(2).ColumnWidth = 49.16666667
$ws.Columns.Item(8).ColumnWidth = 11.16666667

# 3) Overwrite rows 2-11 with the refreshed scrape (new fetch timestamp + new top-10 listing).
#    Several rows in the new scrape have no skill-summary tag, so H is cleared for those.
$ws.Range("A2").Value = "2025-10-23 06:27:52"
$ws.Range("B2").Value = "【急募】APIを利用した診断サイト構築のフリーランス募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5418643"
$ws.Range("G2").Value = 220
$ws.Range("H2").Value = "🔥API ◇サイト"

$ws.Range("A3").Value = "2025-10-23 06:27:52"
$ws.Range("B3").Value = "【急募】施行主向け建築資材配達アプリ開発者を募集します"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5418447"
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = "◆開発 ◇アプリ"

$ws.Range("A4").Value = "2025-10-23 06:27:52"
$ws.Range("B4").Value = "【高品質な恋愛マッチングアプリ制作】エンジニア募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5418455"
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = "◇アプリ"

$ws.Range("A5").Value = "2025-10-23 06:27:52"
$ws.Range("B5").Value = "【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5417544"
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = "◇サイト"

$ws.Range("A6").Value = "2025-10-23 06:27:52"
$ws.Range("B6").Value = "ERPシステムの第三者技術検証・品質評価報告書作成"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5418891"
$ws.Range("G6").Value = 40
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = "2025-10-23 06:27:52"
$ws.Range("B7").Value = "【急募】セッション体験を再現するクローンシステム構築依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5418644"
$ws.Range("G7").Value = 40
$ws.Range("H7").ClearContents()

$ws.Range("A8").Value = "2025-10-23 06:27:52"
$ws.Range("B8").Value = "【急募】既存の予約システムの料金修正を依頼します"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5418759"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

$ws.Range("A9").Value = "2025-10-23 06:27:52"
$ws.Range("B9").Value = "Stable Diffusion LoRA制作依頼 画風再現+キャラLoRA量産テンプレ構築"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5418738"
$ws.Range("G9").Value = 18
$ws.Range("H9").ClearContents()

$ws.Range("A10").Value = "2025-10-23 06:27:52"
$ws.Range("B10").Value = "【メールマーケティング】戦略立案・実行者募集"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5418443"
$ws.Range("G10").Value = 18
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value = "2025-10-23 06:27:52"
$ws.Range("B11").Value = "【急募】HPの微修正をお手伝いしてくれる方募集!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "5,000 円 ~"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5418445"
$ws.Range("G11").Value = 10
$ws.Range("H11").ClearContents()

# 4) Rebuild the URL hyperlinks so only F2:F11 carry links (old F12:F32 links are gone).
#    Re-applying the "Hyperlink" cell style afterwards keeps the cell on the workbook's
#    original Hyperlink style slot instead of a freshly minted duplicate.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5418643", $null, $null, "https://www.lancers.jp/work/detail/5418643")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5418447", $null, $null, "https://www.lancers.jp/work/detail/5418447")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5418455", $null, $null, "https://www.lancers.jp/work/detail/5418455")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5417544", $null, $null, "https://www.lancers.jp/work/detail/5417544")
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5418891", $null, $null, "https://www.lancers.jp/work/detail/5418891")
$ws.Range("F6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5418644", $null, $null, "https://www.lancers.jp/work/detail/5418644")
$ws.Range("F7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5418759", $null, $null, "https://www.lancers.jp/work/detail/5418759")
$ws.Range("F8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5418738", $null, $null, "https://www.lancers.jp/work/detail/5418738")
$ws.Range("F9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5418443", $null, $null, "https://www.lancers.jp/work/detail/5418443")
$ws.Range("F10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5418445", $null, $null, "https://www.lancers.jp/work/detail/5418445")
$ws.Range("F11").Style = "Hyperlink"
